# [ELAB-432] corrections and rewording of some parts^4
#
# Fixes a handful of typos / wording issues in the "Beschreibung" (G) column
# of the Anforderungsliste and moves the active selection, mirroring the
# author's manual editing pass through the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- text corrections -----------------------------------------------------

# Row 4 (Id 2111, "Use-Case-Logging"): "Formular" -> "Formulars"
$ws.Range("G4").Value = 'Tritt im Frontend ein Use-Case auf, soll dieser im Log notiert werden. Beispielsweise soll notiert werden, wenn ein Nutzer das Absenden eines Formulars initiiert.'

# Row 14 (Id 2510, "Session-Replay"): "Schnittstellaufrufe" -> "Schnittstellenaufrufe"
$ws.Range("G14").Value = 'Im Frontend sind Daten zwecks Session-Replay zu erheben, welche u. A. Benutzerinteraktionen, Schnittstellenaufrufe sowie DOM-Manipulationen enthalten.'

# Row 16 (Id 2520, "Übertragung von Session-Replay-Daten"): "Sämtlich" -> "Sämtliche"
$ws.Range("G16").Value = 'Sämtliche im Frontend erfasste Daten zum Session-Replay sind an ein "Session-Replay"-Partnersystem weiterzuleiten.'

# Row 30 (Id 5310, "Manuelle Analyse Tracing"): ", sowie können diese" -> " und können"
$ws.Range("G30").Value = 'Die erfassten Tracingdaten sind für die Nutzer des Systems einsehbar und können gefiltert werden. Die Filtierung erfolgt auf Basis von Eigenschaften der Tracingdaten (wie Name des meldenden Systems).'

# Row 35 (Id 5500, "Partnersystem Session-Replay"): reworded
$ws.Range("G35").Value = 'Es existiert ein "Session-Replay"-Partnersystem, zu dem die Daten zus Session-Replays gesendet werden und welches diese analysiert und speichert.'

# --- view state (selection moved to F38, scrolled so row 22 is at top) ----

$ws.Activate()
$ws.Range("F38").Select()
$excel.ActiveWindow.ScrollRow = 22
